$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update timestamp in A1
$ws.Range("A1").Value = "CreatedAt: 2025-10-29T14:07:54"

# Update numeric data cells
$ws.Range("S4").Value = 99.91
$ws.Range("T4").Value = 109.22
$ws.Range("U4").Value = 97.03
$ws.Range("V4").Value = 73.61
$ws.Range("R6").Value = -4.68
$ws.Range("S6").Value = -4.4
$ws.Range("T6").Value = -5.24
$ws.Range("U6").Value = -4.75
$ws.Range("V6").Value = -3.46
$ws.Range("W6").Value = -3.22
$ws.Range("R9").Value = 102.99
$ws.Range("S9").Value = 102.46
$ws.Range("T9").Value = 111.24
$ws.Range("U9").Value = 99.59
$ws.Range("W9").Value = 69.69
$ws.Range("Y9").Value = 55.98
$ws.Range("R11").Value = -1.34
$ws.Range("S11").Value = -1.84
$ws.Range("T11").Value = -3.23
$ws.Range("U11").Value = -2.19
$ws.Range("Y11").Value = 1.9
$ws.Range("R14").Value = 102.99
$ws.Range("S14").Value = 102.46
$ws.Range("T14").Value = 111.34
$ws.Range("U14").Value = 99.59
$ws.Range("V14").Value = 77.53
$ws.Range("W14").Value = 69.76000000000001
$ws.Range("Y14").Value = 203.52
$ws.Range("Z14").Value = 201.62
$ws.Range("Y15").Value = 147.54
$ws.Range("Z15").Value = 145.64
$ws.Range("R16").Value = -1.34
$ws.Range("S16").Value = -1.84
$ws.Range("T16").Value = -3.12
$ws.Range("U16").Value = -2.19
$ws.Range("V16").Value = 0.47
$ws.Range("Z16").Value = 2.63
$ws.Range("R19").Value = 99.36
$ws.Range("S19").Value = 99.43000000000001
$ws.Range("T19").Value = 108.8
$ws.Range("U19").Value = 96.66
$ws.Range("W19").Value = 65.56999999999999
$ws.Range("R21").Value = -4.97
$ws.Range("S21").Value = -4.87
$ws.Range("T21").Value = -5.66
$ws.Range("U21").Value = -5.12
$ws.Range("W21").Value = -3.34
$ws.Range("R24").Value = 99.36
$ws.Range("S24").Value = 99.43000000000001
$ws.Range("T24").Value = 108.8
$ws.Range("U24").Value = 96.66
$ws.Range("W24").Value = 65.56999999999999
$ws.Range("R26").Value = -4.97
$ws.Range("S26").Value = -4.87
$ws.Range("T26").Value = -5.66
$ws.Range("U26").Value = -5.12
$ws.Range("W26").Value = -3.34
$ws.Range("R29").Value = 98.52
$ws.Range("S29").Value = 98.59
$ws.Range("T29").Value = 107.88
$ws.Range("U29").Value = 95.93000000000001
$ws.Range("W29").Value = 65.2
$ws.Range("Y29").Value = 51.85
$ws.Range("S31").Value = -5.72
$ws.Range("T31").Value = -6.58
$ws.Range("U31").Value = -5.85
$ws.Range("Y31").Value = -2.23
$ws.Range("R34").Value = 105.81
$ws.Range("S34").Value = 105.25
$ws.Range("T34").Value = 113.89
$ws.Range("U34").Value = 102.4
$ws.Range("V34").Value = 80.44
$ws.Range("W34").Value = 72.7
$ws.Range("Y34").Value = 206.26
$ws.Range("Z34").Value = 204.14
$ws.Range("Y35").Value = 147.54
$ws.Range("Z35").Value = 145.64
$ws.Range("R36").Value = 1.48
$ws.Range("S36").Value = 0.95
$ws.Range("T36").Value = -0.57
$ws.Range("U36").Value = 0.61
$ws.Range("V36").Value = 3.38
$ws.Range("Z36").Value = 5.15
$ws.Range("S39").Value = 99.91
$ws.Range("T39").Value = 109.22
$ws.Range("U39").Value = 97.03
$ws.Range("V39").Value = 73.61
$ws.Range("R41").Value = -4.68
$ws.Range("S41").Value = -4.4
$ws.Range("T41").Value = -5.24
$ws.Range("U41").Value = -4.75
$ws.Range("V41").Value = -3.46
$ws.Range("W41").Value = -3.22
$ws.Range("R44").Value = 101.2
$ws.Range("S44").Value = 102.26
$ws.Range("T44").Value = 112.11
$ws.Range("U44").Value = 99.5
$ws.Range("W44").Value = 67.17
$ws.Range("X44").Value = 57.91
$ws.Range("Z44").Value = 52.51
$ws.Range("R46").Value = -3.14
$ws.Range("S46").Value = -2.05
$ws.Range("T46").Value = -2.35
$ws.Range("U46").Value = -2.29
$ws.Range("X46").Value = -1.56
$ws.Range("Z46").Value = -0.84
$ws.Range("R49").Value = 107.67
$ws.Range("S49").Value = 108.99
$ws.Range("T49").Value = 119.73
$ws.Range("U49").Value = 106.58
$ws.Range("W49").Value = 72.78
$ws.Range("Y49").Value = 56.93
$ws.Range("R51").Value = 3.34
$ws.Range("S51").Value = 4.69
$ws.Range("T51").Value = 5.27
$ws.Range("U51").Value = 4.8
$ws.Range("Y51").Value = 2.85
$ws.Range("R54").Value = 102.29
$ws.Range("S54").Value = 102.46
$ws.Range("T54").Value = 112
$ws.Range("U54").Value = 99.69
$ws.Range("W54").Value = 69.69
$ws.Range("X54").Value = 61.25
$ws.Range("R56").Value = -2.05
$ws.Range("T56").Value = -2.46
$ws.Range("U56").Value = -2.09
$ws.Range("X56").Value = 1.78
$ws.Range("R59").Value = 106.79
$ws.Range("S59").Value = 107.53
$ws.Range("T59").Value = 118
$ws.Range("U59").Value = 104.82
$ws.Range("W59").Value = 70.76000000000001
$ws.Range("R61").Value = 2.46
$ws.Range("S61").Value = 3.23
$ws.Range("T61").Value = 3.54
$ws.Range("U61").Value = 3.04
$ws.Range("R64").Value = 108.57
$ws.Range("S64").Value = 109.56
$ws.Range("T64").Value = 120.36
$ws.Range("U64").Value = 106.92
$ws.Range("W64").Value = 72.09
$ws.Range("Y64").Value = 56.39
$ws.Range("R66").Value = 4.23
$ws.Range("S66").Value = 5.26
$ws.Range("T66").Value = 5.9
$ws.Range("U66").Value = 5.13
$ws.Range("W66").Value = 3.17
$ws.Range("Y66").Value = 2.31
$ws.Range("R69").Value = 108.68
$ws.Range("S69").Value = 109.45
$ws.Range("T69").Value = 120.11
$ws.Range("U69").Value = 106.8
$ws.Range("W69").Value = 72.09
$ws.Range("Y69").Value = 56.45
$ws.Range("Z69").Value = 55.75
$ws.Range("R71").Value = 4.35
$ws.Range("S71").Value = 5.14
$ws.Range("T71").Value = 5.64
$ws.Range("U71").Value = 5.02
$ws.Range("W71").Value = 3.17
$ws.Range("Y71").Value = 2.37
$ws.Range("Z71").Value = 2.4
$ws.Range("R74").Value = 107.89
$ws.Range("S74").Value = 108.54
$ws.Range("T74").Value = 119.23
$ws.Range("U74").Value = 105.69
$ws.Range("W74").Value = 71.56999999999999
$ws.Range("R76").Value = 3.56
$ws.Range("S76").Value = 4.23
$ws.Range("T76").Value = 4.77
$ws.Range("U76").Value = 3.91
$ws.Range("W76").Value = 2.65
$ws.Range("R79").Value = 108.29
$ws.Range("S79").Value = 108.8
$ws.Range("T79").Value = 119.61
$ws.Range("U79").Value = 106.27
$ws.Range("Z79").Value = 55.43
$ws.Range("R81").Value = 3.96
$ws.Range("S81").Value = 4.5
$ws.Range("T81").Value = 5.15
$ws.Range("U81").Value = 4.49
$ws.Range("W81").Value = 2.99
$ws.Range("Z81").Value = 2.08
$ws.Range("R84").Value = 99.65000000000001
$ws.Range("S84").Value = 100
$ws.Range("T84").Value = 109.74
$ws.Range("U84").Value = 97.78
$ws.Range("W84").Value = 70.33
$ws.Range("X84").Value = 61.44
$ws.Range("R86").Value = -4.68
$ws.Range("S86").Value = -4.3
$ws.Range("T86").Value = -4.72
$ws.Range("U86").Value = -4.01
$ws.Range("W86").Value = 1.41
$ws.Range("X86").Value = 1.97
$ws.Range("R89").Value = 98.52
$ws.Range("S89").Value = 98.59
$ws.Range("T89").Value = 107.88
$ws.Range("U89").Value = 95.93000000000001
$ws.Range("W89").Value = 65.2
$ws.Range("Y89").Value = 51.85
$ws.Range("S91").Value = -5.72
$ws.Range("T91").Value = -6.58
$ws.Range("U91").Value = -5.85
$ws.Range("Y91").Value = -2.23
